$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "Miss Dina Nasr, Administrator"
$replacement = "Administrator, Miss Dina Nasr"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
